$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cell updates (coin names, links) ---
$ws.Range("B8").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C8").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("B12").Value = "MCDex"
$ws.Range("C12").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"

# --- Numeric-looking text cell updates (price, volume %, hora) ---
# Force text number format per-cell so values stay stored as strings (matching original inline-string cells)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "327.38"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.55%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "2"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "45.18"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.58%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "2"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.341"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-5.57%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "2"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08355"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.20%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "2"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.926"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-6.25%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "2"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9707"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.87%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "2"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.1118"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.70%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "2"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1904"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.10%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "2"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09635"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-4.34%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "2"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.04605"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.02%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "2"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.569"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-26.88%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "2"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1060"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.17%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "2"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001299"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.08%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "2"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006007"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.70%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "2"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.372"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.13%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "2"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.435"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.53%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "2"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.518"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-3.60%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "2"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3357"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.27%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "2"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1388"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.53%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "2"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2598"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.24%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "2"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04170"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.00%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "2"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001252"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-4.15%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "2"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004410"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.34%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "2"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.67%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "2"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-20.32%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "2"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "2"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "2"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "2"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "2"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "2"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "2"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "2"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "2"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "2"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "2"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "2"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02748"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-1.33%"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "2"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05672"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-2.27%"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "2"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007776"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.40%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "2"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1413"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.79%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "2"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006592"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-8.40%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "2"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.01%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "2"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007987"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-5.25%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "2"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3357"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "2"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006980"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.42%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "2"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.04%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "2"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003490"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.31%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "2"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003531"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.96%"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "2"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.04%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "2"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.04%"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "2"
